# Add a second "Logging" sheet (used to persist small bits of config data,
# e.g. a carried-over row/column position) ahead of the existing timesheet
# sheet, and rename the original sheet to "Timesheet".

$wb = $excel.ActiveWorkbook

# Remember the original sheet's name before we start adding/renaming sheets.
$originalName = $wb.Worksheets.Item(1).Name

# Worksheets.Add() inserts the new sheet right before the active sheet -
# exactly where we want "Logging" to land (in front of the timesheet).
$logging = $wb.Worksheets.Add()
$logging.Name = "Logging"

# Re-fetch the original sheet by name (sheet handles obtained before the
# Add() track *position*, not identity, and the insert shifted it along).
$timesheet = $wb.Worksheets.Item($originalName)
$timesheet.Name = "Timesheet"

# Small config table on the new Logging sheet:
#   B1: carryover
#   A2: row      B2: 32
#   A3: column   B3: 10
$logging.Range("B1").Value = "carryover"
$logging.Range("A2").Value = "row"
$logging.Range("B2").Value = 32
$logging.Range("A3").Value = "column"
$logging.Range("B3").Value = 10

# Restore/refresh each sheet's own selection.
[void]$logging.Range("B8").Select()
[void]$timesheet.Range("J32").Select()

# Timesheet is the tab that should be active/visible on open.
$timesheet.Activate()
